$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.010.87'
$ws.Range('E2').Value = '  +3.35%  '
$ws.Range('D3').Value = '1.777.74'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('D4').Value = '''0.9981'
$ws.Range('E4').Value = '  -0.72%  '
$ws.Range('D5').Value = '''332.92'
$ws.Range('E5').Value = '  -1.05%  '
$ws.Range('D6').Value = '''0.9926'
$ws.Range('E6').Value = '  -0.93%  '
$ws.Range('D7').Value = '''0.3816'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').Value = '''0.3422'
$ws.Range('E8').Value = '  +0.53%  '
$ws.Range('D9').Value = '''47.73'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('D10').Value = '''1.145'
$ws.Range('E10').Value = '  -3.28%  '
$ws.Range('D11').Value = '''0.07433'
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '''22.84'
$ws.Range('E12').Value = '  +5.72%  '
$ws.Range('D13').Value = '''0.9917'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('D14').Value = '''6.388'
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '1.775.66'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '''7.128'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('D18').Value = '''0.06632'
$ws.Range('E18').Value = '  +0.02%  '
$ws.Range('D19').Value = '''82.74'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '''0.9963'
$ws.Range('D21').Value = '''17.48'
$ws.Range('E21').Value = '  +1.17%  '
$ws.Range('D22').Value = '''6.431'
$ws.Range('E22').Value = '  -1.47%  '
$ws.Range('D23').Value = '28.034.56'
$ws.Range('E23').Value = '  +3.40%  '
$ws.Range('D24').Value = '''12.12'
$ws.Range('E24').Value = '  -0.56%  '
$ws.Range('D25').Value = '''2.382'
$ws.Range('E25').Value = '  +0.31%  '
$ws.Range('D26').Value = '''1.442'
$ws.Range('E26').Value = '  -0.80%  '
$ws.Range('D27').Value = '''20.85'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('D28').Value = '''2.445'
$ws.Range('E28').Value = '  -2.02%  '
$ws.Range('D29').Value = '''153.98'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('B30').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C30').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D30').Value = '1.977.14'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('B31').Value = 'BitcoinCash'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D31').Value = '''135.01'
$ws.Range('E31').Value = '  +0.96%  '
$ws.Range('D32').Value = '''6.176'
$ws.Range('E32').Value = '  +2.79%  '
$ws.Range('D33').Value = '''3.949'
$ws.Range('E33').Value = '  -0.57%  '
$ws.Range('D34').Value = '''0.08788'
$ws.Range('E34').Value = '  +1.42%  '
$ws.Range('D35').Value = '''12.82'
$ws.Range('E35').Value = '  -1.72%  '
$ws.Range('D36').Value = '''0.02432'
$ws.Range('E36').Value = '  +5.21%  '
$ws.Range('D37').Value = '''0.6887'
$ws.Range('E37').Value = '  +1.07%  '
$ws.Range('D38').Value = '''5.338'
$ws.Range('E38').Value = '  -0.57%  '
$ws.Range('D39').Value = '''0.06357'
$ws.Range('E39').Value = '  +1.41%  '
$ws.Range('D40').Value = '''0.2194'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').Value = '''1.513'
$ws.Range('E41').Value = '  -6.77%  '
$ws.Range('E42').Value = '  +0.54%  '
$ws.Range('D43').Value = '''8.341'
$ws.Range('E43').Value = '  -2.05%  '
$ws.Range('D44').Value = '''14.21'
$ws.Range('E44').Value = '  +0.01%  '
$ws.Range('D45').Value = '''0.9942'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '''0.6319'
$ws.Range('E46').Value = '  -1.43%  '
$ws.Range('D47').Value = '''3.824'
$ws.Range('E47').Value = '  -0.76%  '
$ws.Range('D48').Value = '''132.22'
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').Value = '''2.090'
$ws.Range('E49').Value = '  -1.22%  '
$ws.Range('D50').Value = '''0.07445'
$ws.Range('E50').Value = '  +5.15%  '
$ws.Range('D51').Value = '''78.66'
$ws.Range('E51').Value = '  +0.17%  '
